# "Add test case to excel"
# Populate the QuerySet2 worksheet with additional SingleSynonym test rows
# (rows 3-17), matching the Declaration / Select / Expected Answer / Comment
# columns already established by the header (row 1) and the first sample
# row (row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: StmtStmt, no common synonym
$ws.Range("B3").Value = "stmt s, s1, s2;"
$ws.Range("C3").Value = "Select s with s1.stmt#=s2.stmt#"
$ws.Range("E3").Value = "Select SingleSynonym With StmtStmt NoCommonSynonym HaveResult"

# Row 4: StmtStmt, common synonym
$ws.Range("E4").Value = "Select SingleSynonym With StmtStmt CommonSynonym HaveResult"
$ws.Range("C4").Value = "Select s1 with s1.stmt#=s2.stmt#"
$ws.Range("B4").Value = "stmt s1, s2;"

# Row 5: StmtInt, common synonym
$ws.Range("B5").Value = "stmt s;"
$ws.Range("C5").Value = "Select s with s.stmt#=3"
$ws.Range("E5").Value = "Select SingleSynonym With StmtInt CommonSynonym HaveResult"

# Row 6: StmtInt, no common synonym
$ws.Range("B6").Value = "stmt s1, s2;"
$ws.Range("C6").Value = "Select s1 with s2.stmt#=5"
$ws.Range("E6").Value = "Select SingleSynonym With StmtInt NoCommonSynonym HaveResult"

# Row 7: StmtStr, non compatible type, no result
$ws.Range("B7").Value = "stmt s;"
$ws.Range("C7").Value = "Select s with s.stmt#=`"hello`""
$ws.Range("D7").Value = "none"
$ws.Range("E7").Value = "Select SingleSynonym With StmtStr NonCompatible WithType NoResult"

# Row 8: StmtAssign
$ws.Range("C8").Value = "Select s with s.stmt# = a.stmt#"
$ws.Range("B8").Value = "stmt s; assign a;"
$ws.Range("E8").Value = "Select SingleSynonym With StmtAssign HaveResult"

# Row 9: StmtWhile
$ws.Range("B9").Value = "stmt s; while w;"
$ws.Range("C9").Value = "Select s with s.stmt# = w.stmt#"
$ws.Range("E9").Value = "Select SingleSynonym With StmtWhile HaveResult"

# Row 10: StmtIf
$ws.Range("B10").Value = "stmt s; if f;"
$ws.Range("C10").Value = "Select s with s.stmt# = f.stmt#"
$ws.Range("E10").Value = "Select SingleSynonym With StmtIf HaveResult"

# Row 11: StmtProgLine
$ws.Range("B11").Value = "stmt s; prog_line pl;"
$ws.Range("C11").Value = "Select s with s.stmt# = pl"
$ws.Range("E11").Value = "Select SingleSynonym With StmtProgLine HaveResult"

# Row 12: StmtCall
$ws.Range("B12").Value = "stmt s; call cl;"
$ws.Range("C12").Value = "Select s with s.stmt# = cl.stmt#"
$ws.Range("E12").Value = "Select SingleSynonym With StmtCall HaveResult"

# Row 13: AssignInt
$ws.Range("B13").Value = "assign a; "
$ws.Range("C13").Value = "Select a with a.stmt# = 4"
$ws.Range("E13").Value = "Select SingleSynonym With AssignInt HaveResult"

# Row 14: AssignAssign, lhs/rhs same
$ws.Range("B14").Value = "assign a; "
$ws.Range("C14").Value = "Select a with a.stmt# = a.stmt#"
$ws.Range("E14").Value = "Select SingleSynonym With AssignAssign LhsRhsSame HaveResult"

# Row 15: AssignAssign, common synonym
$ws.Range("B15").Value = "assign a1, a2;"
$ws.Range("C15").Value = "Select a1 with a1.stmt# = a2.stmt#"
$ws.Range("E15").Value = "Select SingleSynonym With AssignAssign CommonSynonym HaveResult"

# Row 16: AssignAssign, no common synonym
$ws.Range("B16").Value = "assign a, a1, a2;"
$ws.Range("C16").Value = "Select a with a1.stmt# = a2.stmt#"
$ws.Range("E16").Value = "Select SingleSynonym With AssignAssign NoCommonSynonym HaveResult"

# Row 17: AssignWhile, no result
$ws.Range("B17").Value = "assign a; while w;"
$ws.Range("C17").Value = "Select a with a.stmt#=w.stmt#"
$ws.Range("E17").Value = "Select SingleSynonym With AssignWhile NoResult"

# Update the view so it reflects scrolling down to the newly-added rows
# and leaves the active selection on C18 (just below the new block).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C18").Select() | Out-Null
